# Update "想去人数" (F column) figures on the "展览", "演出" and "全部类型"
# sheets as produced by the latest data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 87
$ws1.Range("F9").Value  = 8993
$ws1.Range("F10").Value = 823
$ws1.Range("F12").Value = 1153
$ws1.Range("F13").Value = 1021
$ws1.Range("F16").Value = 10
$ws1.Range("F17").Value = 243
$ws1.Range("F18").Value = 307
$ws1.Range("F20").Value = 238
$ws1.Range("F21").Value = 1139

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 87
$ws4.Range("F8").Value  = 6
$ws4.Range("F11").Value = 8993
$ws4.Range("F12").Value = 823
$ws4.Range("F14").Value = 1153
$ws4.Range("F15").Value = 1021
$ws4.Range("F18").Value = 10
$ws4.Range("F19").Value = 243
$ws4.Range("F20").Value = 307
$ws4.Range("F22").Value = 238
$ws4.Range("F23").Value = 1139

$wb.Save()
